$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 31; $r++) {
    $name = $ws.Cells.Item($r, 1).Text

    # Column G: sound reference, e.g. [sound:0001_ദി_01.mp3]
    $ws.Cells.Item($r, 7).Value = "[sound:$name.mp3]"

    # Column H: image tag, e.g. <img src="0001_ദി_01.jpg">
    $ws.Cells.Item($r, 8).Value = '<img src="' + $name + '.jpg">'
}
